$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'288.91"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-0.37%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'31.03"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'0.51%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'4.932"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-0.15%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.07380"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'2.51%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'2.226"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'24.55%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'7.698"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'0.42%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'3.733"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'0.28%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.9088"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'1.46%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.08753"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'13.34%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.1685"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'1.74%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.08290"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'3.06%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.03113"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'2.54%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.09942"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'-0.73%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.001496"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-0.20%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.005703"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-1.46%"
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = "'0.51%"
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Value = "'0.33%"
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'0.49%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.1298"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'0.06%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'3.828"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'-5.31%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.2122"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'1.03%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.04552"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'0.64%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.001211"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'-0.31%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.004147"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'3.54%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.0001301"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'4.01%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'0.0003396"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'-95.49%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D39").Value = "'0.01581"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'-0.67%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.04474"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'1.98%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.007331"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'-0.15%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.009558"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'24.56%"
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'1.39%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.002141"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'6.34%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.008331"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-9.45%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006106"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'3.25%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.00000000750"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'0.03%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'2.260"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'0.64%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.002001"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'-33.31%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.00002101"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'0.03%"
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.0002001"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'0.03%"
$ws.Range("E51").Style = "Normal"
